$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Very noisy"
$ws.Range("H7").Value = "Very noisy"
$ws.Range("H10").Value = "kinda noisy"
$ws.Range("H14").Value = "kinda noisy"
$ws.Range("H15").Value = "Very noisy"
$ws.Range("H16").Value = "Very noisy"
$ws.Range("H20").Value = "Very noisy"
$ws.Range("H21").Value = "kinda noisy"
$ws.Range("H25").Value = "Very noisy"
$ws.Range("H26").Value = "kinda noisy"
$ws.Range("H29").Value = "kinda noisy"
$ws.Range("H31").Value = "kinda noisy"
$ws.Range("H34").Value = "very noisy"
$ws.Range("H35").Value = "kinda noisy"
$ws.Range("H38").Value = "kinda noisy"
$ws.Range("H40").Value = "very noisy"
$ws.Range("H49").Value = "kinda noisy"
$ws.Range("H58").Value = "Very noisy"
$ws.Range("H62").Value = "kinda noisy"
$ws.Range("G68").Value = "kinda noisy"
$ws.Range("G69").Value = "kinda noisy"

$ws.Range("G71").Select()
